# Auto-update the "剩余" (remaining days) / "开始时间" (start date) columns to
# reflect the passage of time to a new reference ("today") date of 2025-11-25.
#
# Business rule recovered from the data (D = total days, E = remaining days,
# F = start date as YYYYMMDD):
#     end_date = F + D days
#     E        = end_date - today   (in whole days)
# When that would make E <= 0 (the period already expired as of the new
# "today"), the record is renewed: F resets to "today" and E resets back to
# the full D (i.e. a fresh period starting today).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = Get-Date -Year 2025 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0
$todayOA = [math]::Floor($today.ToOADate())
$todaySerial = 20251125

function Get-OADateFromYyyymmdd($v) {
    # Parses an integer formatted as YYYYMMDD into a truncated OADate
    # (day-count) value. Returns $null if the value isn't a well-formed
    # 8-digit date (e.g. the corrupt "202510929" in row 36), so that row
    # can be left untouched, same as in the source diff.
    if ($v -eq $null) { return $null }
    $s = [string]([long]$v)
    if ($s.Length -ne 8) { return $null }
    $y = [int]$s.Substring(0,4)
    $mo = [int]$s.Substring(4,2)
    $da = [int]$s.Substring(6,2)
    if ($mo -lt 1 -or $mo -gt 12) { return $null }
    if ($da -lt 1 -or $da -gt 31) { return $null }
    try {
        $dt = Get-Date -Year $y -Month $mo -Day $da -Hour 0 -Minute 0 -Second 0
        return [math]::Floor($dt.ToOADate())
    } catch {
        return $null
    }
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $fVal -eq $null) { continue }

    $fOA = Get-OADateFromYyyymmdd $fVal
    if ($fOA -eq $null) { continue }

    $endOA = $fOA + [double]$dVal
    $remaining = $endOA - $todayOA

    if ($remaining -le 0) {
        # Expired as of the new reference date: renew the period starting today.
        $eCell.Value = [double]$dVal
        $fCell.Value = $todaySerial
    } else {
        $eCell.Value = $remaining
    }
}
